$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table refresh (GitHub Actions data pull).
# Two rows (18/19, Uniswap <-> WrappedEther) swapped rank position; all other
# rows keep their coin/link but get refreshed Price (D) / Volume 1h (E) values.

$ws.Range("D2").Value = "70.816.10"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "3.555.25"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'582.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "'187.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("D8").Value = "3.545.31"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "'0.218"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.49%  "
$ws.Range("D11").Value = "'0.650"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "'54.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("E13").Value = "  +4.80%  "
$ws.Range("D14").Value = "'9.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "4.120.33"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "70.848.72"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "'19.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.558.31"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'12.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("D20").Value = "'576.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.09%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("E23").Value = "  -5.02%  "
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").Value = "'4.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D26").Value = "'94.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'11.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'9.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "'32.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "'63.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  +22.79%  "
$ws.Range("E36").Value = "  +7.77%  "
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").Value = "'530.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").Value = "'38.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "0.0₃0800"
$ws.Range("E40").Value = "  +4.13%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "3.620.55"
$ws.Range("E42").Value = "  +9.36%  "
$ws.Range("E43").Value = "  +3.40%  "
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").Value = "'0.0467"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.57%  "
$ws.Range("D46").Value = "'3.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("E51").Value = "  +6.02%  "
